$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 7.2
$ws.Range("H5").Value = 'От 400 ₽'
$ws.Range("B9").Value = 'Ревизор'
$ws.Range("C9").Value = '26 июня'
$ws.Range("D9").Value = '12:00'
$ws.Range("E9").Value = 'Комедия'
$ws.Range("G9").Value = 8.1
$ws.Range("H9").Value = 'От 400 ₽'
$ws.Range("I9").Value = 'Театр им. Ленсовета'
$ws.Range("J9").Value = 'Презабавный утренник по пьесе Гоголя, который отчего-то играют вечером'
$ws.Range("K9").Value = 'https://www.afisha.ru/performance/78711/'
$ws.Range("B10").Value = 'Лето одного года'
$ws.Range("C10").Value = '8 апреля'
$ws.Range("F10").Value = '12+'
$ws.Range("G10").Value = 8.9
$ws.Range("H10").Value = 'От 4000 ₽'
$ws.Range("J10").Value = 'Алиса Фрейндлих и Олег Басилашвили с блеском и мужеством играют трагикомедию о старости и об уходе вообще'
$ws.Range("K10").Value = 'https://www.afisha.ru/performance/82546/'
$ws.Range("B11").Value = 'Топливо'
$ws.Range("C11").Value = '31 марта'
$ws.Range("D11").Value = '20:00'
$ws.Range("E11").Value = 'Драма'
$ws.Range("G11").Value = 8.4
$ws.Range("H11").Value = 'От 1000 ₽'
$ws.Range("I11").Value = 'Скороход'
$ws.Range("J11").Value = 'Увлекательный байопик про российского физика'
$ws.Range("K11").Value = 'https://www.afisha.ru/performance/104731/'
$ws.Range("B12").Value = 'Гроза'
$ws.Range("C12").Value = '28 апреля'
$ws.Range("E12").Value = 'Драма'
$ws.Range("F12").Value = '16+'
$ws.Range("G12").Value = 8
$ws.Range("H12").Value = 'От 300 ₽'
$ws.Range("I12").Value = 'БДТ'
$ws.Range("J12").Value = 'Большое шаманство Андрея Могучего по Островскому'
$ws.Range("K12").Value = 'https://www.afisha.ru/performance/117877/'
$ws.Range("B13").Value = 'С Чарльзом Буковски за барной стойкой'
$ws.Range("C13").Value = 'Сегодня'
$ws.Range("E13").Value = 'Перформанс'
$ws.Range("F13").Value = '18+'
$ws.Range("G13").Value = 5.2
$ws.Range("H13").Value = 'От 900 ₽'
$ws.Range("I13").Value = 'Fiddler''s Green'
$ws.Range("J13").Value = 'Спектакль для одного зрителя в наушниках за барной стойкой'
$ws.Range("K13").Value = 'https://www.afisha.ru/performance/191308/'
$ws.Range("B14").Value = 'Безымянная звезда'
$ws.Range("C14").Value = '7 апреля'
$ws.Range("F14").Value = '14+'
$ws.Range("H14").Value = 'От 300 ₽'
$ws.Range("I14").Value = 'Театр им. Комиссаржевской'
$ws.Range("J14").Value = 'Первая любовь гения как неисчерпаемый источник поэзии'
$ws.Range("K14").Value = 'https://www.afisha.ru/performance/67513/'
$ws.Range("B15").Value = 'Дети солнца'
$ws.Range("C15").Value = '23 апреля'
$ws.Range("F15").Value = '16+'
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 'От 300 ₽'
$ws.Range("I15").Value = 'Александринский театр'
$ws.Range("J15").Value = 'Размышления об идеальном человеке'
$ws.Range("K15").Value = 'https://www.afisha.ru/performance/233405/'
$ws.Range("B16").Value = 'Что делать'
$ws.Range("C16").Value = '29 марта'
$ws.Range("G16").Value = 6.2
$ws.Range("H16").Value = 'От 500 ₽'
$ws.Range("I16").Value = 'БДТ'
$ws.Range("J16").Value = 'Разговорная антиутопия Андрея Могучего по Чернышевскому'
$ws.Range("K16").Value = 'https://www.afisha.ru/performance/101454/'
$ws.Range("B18").Value = 'Город. Женитьба. Гоголь'
$ws.Range("C18").Value = '29 апреля'
$ws.Range("D18").Value = '19:00'
$ws.Range("E18").Value = 'Драма'
$ws.Range("F18").Value = '12+'
$ws.Range("G18").Value = 8.199999999999999
$ws.Range("H18").Value = 'От 600 ₽'
$ws.Range("I18").Value = 'Театр им. Ленсовета'
$ws.Range("J18").Value = 'Герои Гоголя растворяются в современном Петербурге'
$ws.Range("K18").Value = 'https://www.afisha.ru/performance/103390/'
$ws.Range("B19").Value = 'Когда я снова стану маленьким'
$ws.Range("C19").Value = '24 апреля'
$ws.Range("D19").Value = '13:00'
$ws.Range("E19").Value = 'Кукольный спектакль'
$ws.Range("F19").Value = '6+'
$ws.Range("G19").Value = 8.300000000000001
$ws.Range("H19").Value = 'Билеты'
$ws.Range("J19").Value = 'Энциклопедия подростковых грез и слез, упакованная в 15 чемоданов'
$ws.Range("K19").Value = 'https://www.afisha.ru/performance/97786/'
$ws.Range("B20").Value = 'Губернатор'
$ws.Range("C20").Value = '9 апреля'
$ws.Range("E20").Value = 'Драма'
$ws.Range("F20").Value = '18+'
$ws.Range("G20").Value = 8.6
$ws.Range("H20").Value = 'От 500 ₽'
$ws.Range("I20").Value = 'БДТ'
$ws.Range("J20").Value = 'Зрелищный спектакль Андрея Могучего про механику тоталитаризма'
$ws.Range("K20").Value = 'https://www.afisha.ru/performance/190192/'
$ws.Range("B21").Value = 'Все мы прекрасные люди'
$ws.Range("C21").Value = '19 апреля'
$ws.Range("E21").Value = 'Комедия'
$ws.Range("F21").Value = '16+'
$ws.Range("G21").Value = 4.5
$ws.Range("H21").Value = 'От 600 ₽'
$ws.Range("I21").Value = 'Театр им. Ленсовета'
$ws.Range("J21").Value = 'Внезапная страсть героини Анны Ковальчук сокрушает ее саму и все вокруг'
$ws.Range("K21").Value = 'https://www.afisha.ru/performance/96366/'
$ws.Range("B23").Value = 'Подыскиваю жену, недорого!'
$ws.Range("C23").Value = '16 апреля'
$ws.Range("E23").Value = 'Комедия'
$ws.Range("G23").Value = 5.8
$ws.Range("H23").Value = 'От 1200 ₽'
$ws.Range("I23").Value = 'ДК им. Ленсовета'
$ws.Range("J23").Value = 'Антрепризный спектакль о превратностях любви'
$ws.Range("K23").Value = 'https://www.afisha.ru/performance/85589/'
